# Update Top10 league_soccer and NBA
# Applies the changes described in the commit to planilha_NBA.xlsx

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Jogos de Hoje": clear out the stale list of today's games,
# keep only the header row, and shrink column A back down.
# ---------------------------------------------------------------------
$wsJogos = $wb.Worksheets.Item("Jogos de Hoje")
$wsJogos.Range("A2:C14").EntireRow.Delete() | Out-Null
# 13.14 round-trips through the engine's width<->pixel conversion to
# land on an exact stored width of 14 (matching the target column width).
$wsJogos.Columns.Item(1).ColumnWidth = 13.14

# ---------------------------------------------------------------------
# Sheet "Arremessos de 3 Pontos"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Arremessos de 3 Pontos")
$ws.Range("D3").Value = "'3.9"
$ws.Range("A6").Value = 5
$ws.Range("D6").Value = "'3.3"

# ---------------------------------------------------------------------
# Sheet "Assistências"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Assistências")
$ws.Range("D3").Value = "'10.9"
$ws.Range("D4").Value = "'9.5"
$ws.Range("D5").Value = "'8.9"
$ws.Range("D6").Value = "'8.4"

# ---------------------------------------------------------------------
# Sheet "Pontos"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Pontos")
$ws.Range("D3").Value = "'34.2"
$ws.Range("D5").Value = "'30.8"
$ws.Range("B6").Value = "Donovan Mitchell"
$ws.Range("C6").Value = "CLE"
$ws.Range("D6").Value = "'28.4"

# ---------------------------------------------------------------------
# Sheet "Rebotes"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Rebotes")
$ws.Range("B4").Value = "Anthony Davis"
$ws.Range("C4").Value = "LAL"
$ws.Range("A5").Value = 4
$ws.Range("D5").Value = "'12.0"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Jalen Duren"
$ws.Range("C6").Value = "DET"
$ws.Range("D6").Value = "'12.0"

# ---------------------------------------------------------------------
# Sheet "Roubos"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Roubos")
$ws.Range("A4").Value = 2
$ws.Range("D4").Value = "'1.9"

# ---------------------------------------------------------------------
# Sheet "Tocos"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tocos")
$ws.Range("D2").Value = "'3.2"
$ws.Range("D6").Value = "'2.5"

# ---------------------------------------------------------------------
# Sheet "Arremessos %"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Arremessos %")
$ws.Range("B3").Value = "Oklahoma City Thunder"
$ws.Range("C3").Value = "'50.0"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Phoenix Suns"
$ws.Range("C5").Value = "'49.4"
$ws.Range("B6").Value = "Los Angeles Lakers"
$ws.Range("C6").Value = "'49.4"

# ---------------------------------------------------------------------
# Sheet "Diferencial de Pontos"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Diferencial de Pontos")
$ws.Range("C2").Value = "'+10.1"
$ws.Range("B3").Value = "Minnesota Timberwolves"
$ws.Range("C3").Value = "'+7.2"
$ws.Range("B4").Value = "Oklahoma City Thunder"
$ws.Range("C4").Value = "'+7.1"
$ws.Range("B5").Value = "Cleveland Cavaliers"
$ws.Range("C5").Value = "'+5.6"
$ws.Range("B6").Value = "LA Clippers"
$ws.Range("C6").Value = "'+5.5"

# ---------------------------------------------------------------------
# Sheet "Pontos1"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Pontos1")
$ws.Range("C2").Value = "'123.7"
$ws.Range("C3").Value = "'122.0"
$ws.Range("C4").Value = "'121.3"
$ws.Range("C5").Value = "'120.8"
$ws.Range("C6").Value = "'120.7"

# ---------------------------------------------------------------------
# Sheet "Pontos Permitidos"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Pontos Permitidos")
$ws.Range("C2").Value = "'106.7"
$ws.Range("C3").Value = "'109.2"
$ws.Range("B5").Value = "Miami Heat"
$ws.Range("C5").Value = "'110.4"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Boston Celtics"
$ws.Range("C6").Value = "'110.6"

# ---------------------------------------------------------------------
# Sheet "Rebotes1"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Rebotes1")
$ws.Range("C2").Value = "'47.3"
$ws.Range("A4").Value = 3
$ws.Range("C4").Value = "'46.6"
$ws.Range("C5").Value = "'45.9"
$ws.Range("C6").Value = "'45.7"

# ---------------------------------------------------------------------
# Sheet "Tocos1"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tocos1")
$ws.Range("B2").Value = "Oklahoma City Thunder"
$ws.Range("C2").Value = "'6.7"
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Boston Celtics"
$ws.Range("C3").Value = "'6.5"
$ws.Range("C4").Value = "'6.2"
$ws.Range("B5").Value = "Phoenix Suns"
$ws.Range("C5").Value = "'6.1"
$ws.Range("B6").Value = "San Antonio Spurs"
$ws.Range("C6").Value = "'6.1"

# ---------------------------------------------------------------------
# Sheet "Doubledouble"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Doubledouble")
$ws.Range("D2").Value = 50
$ws.Range("D4").Value = 43
$ws.Range("D5").Value = 42
$ws.Range("D6").Value = 38
